$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column B (N column), shifting everything right.
$ws.Range("B:C").Insert()

# Move the SMD columns (now shifted to G and H after the insert) into the new B and C columns.
$ws.Range("G:H").Cut($ws.Range("B1"))

# Remove the now-empty columns left behind by the cut so the sheet shrinks back down.
$ws.Range("G:H").Delete()

# Update the header labels for the relocated columns.
$ws.Range("B1").Value = "DSM.scale_trim"
$ws.Range("C1").Value = "DSM.scale_trim_round"
